# Add "POWER Data Management Hub" entry to the "slides" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("slides")
$ws.Activate()

$ws.Cells.Item(7, 1).Value = "POWER Data Management Hub"
$ws.Cells.Item(7, 2).Value = "Slides from hub presenters"

$ws.Hyperlinks.Add($ws.Cells.Item(7, 2), "https://rdmkit.elixir-europe.org/", $null, $null, $null)

$ws.Range("B12").Select()
